# Applies the edits described in the commit:
#  1. "Rollen (Admin, Kunde, Mitarbeiter,..)" - the sentence used to be split
#     across several runs around a grammar-check annotation on "Mitarbeiter".
#     Re-asserting the full sentence text makes Word coalesce the runs and
#     drop the stale <w:proofErr/> markers.
#  2. "Kunden (... Instazugriff)" - same idea, "Instazugriff" was wrapped in
#     spell-check markers; re-typing the whole phrase merges the runs and
#     clears the spelling annotation.
#  3. A new bullet "TEST TEST TEST" is appended right after the existing
#     "Arbeitszeiten aktualisieren" bullet, at the same outline level.

$d = $word.ActiveDocument

# --- Edit 1: "Rollen (...)" paragraph -------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(
    "Rollen (Admin, Kunde, Mitarbeiter,..)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Rollen (Admin, Kunde, Mitarbeiter,..)", 2) | Out-Null

# --- Edit 2: "Kunden (...)" paragraph --------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "Kunden (Termine reservieren, offene Termine ansehen, Infos holen, Instazugriff)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Kunden (Termine reservieren, offene Termine ansehen, Infos holen, Instazugriff)",
    2) | Out-Null

# --- Edit 3: append a new bullet after "Arbeitszeiten aktualisieren" -------
# Locate the "Arbeitszeiten aktualisieren" bullet (falls back to the very
# last paragraph of the body if, for some reason, it can't be found) and add
# a new sibling bullet right after it, inheriting its list formatting.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text -replace "[\x07\x0d]", ""
    if ($t -eq "Arbeitszeiten aktualisieren") {
        $target = $p
    }
}
if ($target -eq $null) {
    $target = $d.Paragraphs($d.Paragraphs.Count)
}

$endRange = $target.Range
$endRange.Collapse(0)              # wdCollapseEnd
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertBefore("TEST TEST TEST")

Write-Host "Edits applied"
